# This script applies the edits described by the diff:
#  - 9 "surplus" multiplier cells (K16, K17, K20, K23, K27, K30, K33, K34, K38)
#    are changed from 1.0565 to 1 (the commit message calls this "fixed
#    surplus number"). K25/K26/etc. were already 1 and are left untouched.
#  - The active selection on the sheet moves from A38 to K39.
#  - The shared-formula "ref" span recorded for the E25/G25 shared-formula
#    group is tightened from E25:E31 / G25:G31 down to E25:E26 / G25:G26 to
#    correctly reflect that only E25:E26/G25:G26 actually share that
#    formula (E27 already carries its own standalone copy of the same
#    formula, and E28:E31/G28:G31 are a separate shared group). We rewrite
#    the formulas for that exact two-cell block as a single range-formula
#    assignment, which is how this object model derives/records a shared
#    formula's ref span.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Fixed surplus number": 1.0565 -> 1 --------------------------------
$ws.Range("K16").Value = 1
$ws.Range("K17").Value = 1
$ws.Range("K20").Value = 1
$ws.Range("K23").Value = 1
$ws.Range("K27").Value = 1
$ws.Range("K30").Value = 1
$ws.Range("K33").Value = 1
$ws.Range("K34").Value = 1
$ws.Range("K38").Value = 1

# --- Support for longer quotes: tighten the E25/G25 shared formula group
$ws.Range("E25:E26").Formula = "=ROUND(J25*0.7*K25/0.6*1.1*exchange, 0)"
$ws.Range("G25:G26").Formula = "=E25*F25"

# --- Update the selected cell on the sheet ------------------------------
$ws.Range("K39").Select()
